$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet 1")

# --- Row 5 additions (new scan results appended to an existing row) ---
$ws.Range("X5").Value = 0.11999500000000296
$ws.Range("Y5").Value = "Up"

# --- New row 6 (full new scan result row) ---
# Copy formatting from row 5 first so date/percent number formats are reused
# (same style indices) rather than minting brand-new styles.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("S5").Copy($ws.Range("S6"))
$ws.Range("T5").Copy($ws.Range("T6"))

$ws.Range("A6").Value = 42647.886967592596
$ws.Range("B6").Value = -3
$ws.Range("C6").Value = "Neutral"
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = 15158
$ws.Range("F6").Value = 2686
$ws.Range("G6").Value = 66
$ws.Range("H6").Value = 33
$ws.Range("I6").Value = 75
$ws.Range("J6").Value = 23
$ws.Range("K6").Value = 29731
$ws.Range("L6").Value = 330
$ws.Range("M6").Value = 168
$ws.Range("N6").Value = 84
$ws.Range("O6").Value = 26
$ws.Range("P6").Value = "Noun"
$ws.Range("Q6").Value = 42.459412013272512
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = -0.0112
$ws.Range("T6").Value = -0.0367
$ws.Range("U6").Value = 14.56
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = -2
